# LDLC smartphone price-tracking workbook:
# a new hourly snapshot column is inserted right before the existing
# "nom" / "url_produit" columns (previously AY/AZ, now shifted to AZ/BA).
# The newly inserted column (AY) gets the timestamp of the new scrape in
# row 1, and for every product row that already had data in the last
# existing snapshot column (AX), the same price is duplicated into the
# new AY column (rows without previous data stay empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AY; this automatically shifts the old AY ("nom")
# to AZ and the old AZ ("url_produit") to BA, and grows the used range
# from A1:AZ206 to A1:BA206.
$ws.Columns("AY:AY").Insert()

# Header: timestamp of the new snapshot column.
$ws.Range("AY1").Value = "2026-01-29 22:15:58"

# Carry forward the last known price (column AX) into the new column AY
# for every product row that already had a price tracked through AX.
for ($row = 2; $row -le 80; $row++) {
    $lastPrice = $ws.Cells.Item($row, 50).Value2
    $ws.Cells.Item($row, 51).Value = $lastPrice
}
